$d = $word.ActiveDocument

# Locate the run of text that starts with "in R. My name is Clay Ford..."
# by restricting a Range to the first paragraph, starting right after
# the "Linear Modeling " run, so the Find/Replace below only affects
# that single run and doesn't merge it with neighboring runs.
$full = $d.Content.Text
$startIdx = $full.IndexOf("in R. My name is Clay Ford")
$p1End = $d.Paragraphs.Item(1).Range.End
$r = $d.Range($startIdx, $p1End)

$find = $r.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Statistical Research Consultant ", $true, $false, $false, $false, $false, $true, 1, $false, "statistician ", 2)
